# Weekly refresh: Fecha/Volumen/Precio fields reshuffled across rows 2-27
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: Row, D (Fecha serial), J (Volumen), K (Precio minimo),
#          L (Precio maximo), M (Precio promedio ponderado), P (Precio $/Kg)
$data = @(
    ,@(2, 44838, 10, 20000, 20000, 20000, 1333)
    ,@(3, 44525, 40, 8000, 8000, 8000, 533)
    ,@(4, 44827, 20, 20000, 20000, 20000, 1333)
    ,@(5, 44830, 25, 12000, 12000, 12000, 800)
    ,@(6, 44819, 100, 20000, 20000, 20000, 1333)
    ,@(7, 44749, 50, 20000, 20000, 20000, 1333)
    ,@(8, 44757, 30, 20000, 20000, 20000, 1333)
    ,@(9, 44825, 30, 20000, 20000, 20000, 1333)
    ,@(10, 45134, 5, 20000, 20000, 20000, 1333)
    ,@(11, 44813, 20, 20000, 20000, 20000, 1333)
    ,@(12, 44841, 20, 16000, 16000, 16000, 1067)
    ,@(13, 44839, 80, 16000, 16000, 16000, 1067)
    ,@(14, 44826, 50, 20000, 20000, 20000, 1333)
    ,@(15, 44812, 80, 20000, 20000, 20000, 1333)
    ,@(16, 44845, 20, 16000, 16000, 16000, 1067)
    ,@(17, 44755, 50, 20000, 20000, 20000, 1333)
    ,@(18, 44776, 80, 20000, 20000, 20000, 1333)
    ,@(19, 44769, 50, 20000, 20000, 20000, 1333)
    ,@(20, 44811, 30, 20000, 20000, 20000, 1333)
    ,@(21, 44824, 20, 20000, 20000, 20000, 1333)
    ,@(22, 44508, 40, 10000, 10000, 10000, 667)
    ,@(23, 44837, 80, 16000, 16000, 16000, 1067)
    ,@(24, 44518, 50, 10000, 10000, 10000, 667)
    ,@(25, 44771, 40, 20000, 20000, 20000, 1333)
    ,@(26, 44756, 80, 20000, 20000, 20000, 1333)
    ,@(27, 44767, 50, 20000, 20000, 20000, 1333)
)

foreach ($row in $data) {
    $r = $row[0]
    $d = $row[1]
    $j = $row[2]
    $k = $row[3]
    $l = $row[4]
    $m = $row[5]
    $p = $row[6]
    $ws.Cells.Item($r, 4).Value2  = $d   # D: Fecha
    $ws.Cells.Item($r, 10).Value  = $j   # J: Volumen
    $ws.Cells.Item($r, 11).Value  = $k   # K: Precio minimo
    $ws.Cells.Item($r, 12).Value  = $l   # L: Precio maximo
    $ws.Cells.Item($r, 13).Value  = $m   # M: Precio promedio ponderado
    $ws.Cells.Item($r, 16).Value  = $p   # P: Precio $/Kg
}
